$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeRef, $value) {
    $rng = $ws.Range($rangeRef)
    # Force the cell to stay text even if the value looks like a number or a date,
    # then restore the "Normal" style so no stray formatting is left behind.
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 ------------------------------------------------------------
Set-TextValue "A2" "Sell"
Set-TextValue "G2" "MCX"
Set-TextValue "H2" "CRUDEOIL"
Set-TextValue "I2" "2022-08-19"
Set-TextValue "J2" "CRUDEOIL22AUGFUT"
Set-TextValue "P2" "2"
Set-TextValue "Q2" "Percentage"
Set-TextValue "R2" "0.2"
Set-TextValue "S2" "Percentage"
Set-TextValue "T2" "25"

# Row 3 ------------------------------------------------------------
Set-TextValue "A3" "Buy"
Set-TextValue "G3" "MCX"
Set-TextValue "H3" "CRUDEOIL"
Set-TextValue "I3" "2022-08-19"
Set-TextValue "J3" "CRUDEOIL22AUGFUT"
Set-TextValue "O3" "Percentage"
Set-TextValue "P3" "2"
Set-TextValue "Q3" "Percentage"
Set-TextValue "R3" "0.2"
Set-TextValue "S3" "Percentage"
Set-TextValue "T3" "25"
